$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update D5 value: "xlsx.v2" -> "xlsx.v2\n" (append literal backslash-n)
$ws.Range("D5").Value = "xlsx.v2\n"

# Update selection to G7
$ws.Range("G7").Select()
